# Update "合肥-漫展信息" workbook to the latest scraped snapshot.
# Sheet 1: 展览 (Exhibitions)   -> rows shift up by 2 (oldest 2 events dropped)
# Sheet 2: 演出 (Performances)  -> only want-to-go counts (F) bumped for two rows
# Sheet 3: 本地生活 (Local life) -> unaffected (header only)
# Sheet 4: 全部类型 (All types)  -> rows shift up by 2 (oldest 2 events dropped)

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        [int]$idx,
        [string]$date,
        [string]$name,
        [string]$place,
        [string]$timeRange,
        $want,
        $price,
        [string]$link,
        [string]$cover
    )
    $ws.Cells.Item($row, 1).Value = $idx
    # Column B holds a plain "YYYY-MM-DD" label, not a real date serial, in
    # the source data - force text formatting first so Excel doesn't
    # auto-convert the literal string into a date value on input.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $date
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $place
    $ws.Cells.Item($row, 5).Value = $timeRange
    $ws.Cells.Item($row, 6).Value = $want
    $ws.Cells.Item($row, 7).Value = $price
    $ws.Cells.Item($row, 8).Value = $link
    $ws.Cells.Item($row, 9).Value = $cover
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Set-Row $ws1 2 1 "2024-08-24" "合肥·微光mini动漫派对（免费展）" `
    "山林路与山水路交叉路口往东北约70米 伟星星悦广场(肥东店)" `
    "2024.08.24 13:00-08.25 19:00" 105 58 `
    "https://show.bilibili.com/platform/detail.html?id=90625" `
    "//i0.hdslb.com/bfs/openplatform/202408/t7kq4X7h1723471019389.jpeg"

Set-Row $ws1 3 2 "2024-09-07" "合肥·国乙only宇宙心动（含夜场）" `
    "文忠路1865号 赫拉诺言艺术中心" `
    "2024.09.07 10:00-09.07 21:00" 413 48 `
    "https://show.bilibili.com/platform/detail.html?id=89803" `
    "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

Set-Row $ws1 4 3 "2024-09-15" "合肥·书香璃樱动漫游戏嘉年华" `
    "阜阳北路与金海路交口 格律诗宴会大酒店(北城店)" `
    "2024.09.15 10:00-09.15 17:00" 11 50 `
    "https://show.bilibili.com/platform/detail.html?id=90735" `
    "//i2.hdslb.com/bfs/openplatform/202408/7alsu0yg1723110506313.jpeg"

Set-Row $ws1 5 4 "2024-09-15" "合肥·曙光次元动漫游戏嘉年华" `
    "田埠西路199号 吉祥如意宴会楼蜀山店" `
    "2024.09.15 10:00-09.15 17:00" 8 50 `
    "https://show.bilibili.com/platform/detail.html?id=90733" `
    "//i1.hdslb.com/bfs/openplatform/202408/bNZ6vKL01723113544322.jpeg"

Set-Row $ws1 6 5 "2024-09-16" "肥西·星域动漫游戏嘉年华" `
    "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)" `
    "2024.09.16 10:00-09.16 17:00" 24 45 `
    "https://show.bilibili.com/platform/detail.html?id=90489" `
    "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

# Rows 7-8 no longer exist; delete the now-unused trailing rows entirely.
$ws1.Rows("7:8").Delete()

# ---------------------------------------------------------------------------
# Sheet 2: 演出  (only the "want-to-go" counts move, everything else the same)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 65
$ws2.Cells.Item(3, 6).Value = 25

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 -> unchanged
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Set-Row $ws4 2 1 "2024-08-24" "合肥·微光mini动漫派对（免费展）" `
    "山林路与山水路交叉路口往东北约70米 伟星星悦广场(肥东店)" `
    "2024.08.24 13:00-08.25 19:00" 105 58 `
    "https://show.bilibili.com/platform/detail.html?id=90625" `
    "//i0.hdslb.com/bfs/openplatform/202408/t7kq4X7h1723471019389.jpeg"

Set-Row $ws4 3 2 "2024-08-25" "合肥·CrossingX意次元｜乐队番ONLY同人" `
    "国祯广场B-1楼 背影骑士LIVEHOUSE" `
    "2024.08.25 13:30-08.25 16:00" 65 38 `
    "https://show.bilibili.com/platform/detail.html?id=90032" `
    "//i2.hdslb.com/bfs/openplatform/202407/GYPAnumr1721896597593.jpeg"

Set-Row $ws4 4 3 "2024-09-07" "合肥·国乙only宇宙心动（含夜场）" `
    "文忠路1865号 赫拉诺言艺术中心" `
    "2024.09.07 10:00-09.07 21:00" 413 48 `
    "https://show.bilibili.com/platform/detail.html?id=89803" `
    "//i1.hdslb.com/bfs/openplatform/202407/w5hQDj821721564303601.jpeg"

Set-Row $ws4 5 4 "2024-09-15" "合肥·书香璃樱动漫游戏嘉年华" `
    "阜阳北路与金海路交口 格律诗宴会大酒店(北城店)" `
    "2024.09.15 10:00-09.15 17:00" 11 50 `
    "https://show.bilibili.com/platform/detail.html?id=90735" `
    "//i2.hdslb.com/bfs/openplatform/202408/7alsu0yg1723110506313.jpeg"

Set-Row $ws4 6 5 "2024-09-15" "合肥·曙光次元动漫游戏嘉年华" `
    "田埠西路199号 吉祥如意宴会楼蜀山店" `
    "2024.09.15 10:00-09.15 17:00" 8 50 `
    "https://show.bilibili.com/platform/detail.html?id=90733" `
    "//i1.hdslb.com/bfs/openplatform/202408/bNZ6vKL01723113544322.jpeg"

Set-Row $ws4 7 6 "2024-09-16" "肥西·星域动漫游戏嘉年华" `
    "金寨路与云谷路交口金云国际9号楼商(邮政银行旁边) 吉祥如意大酒店(肥西店)" `
    "2024.09.16 10:00-09.16 17:00" 24 45 `
    "https://show.bilibili.com/platform/detail.html?id=90489" `
    "//i2.hdslb.com/bfs/openplatform/202408/6xk6G8E71722525186252.jpeg"

Set-Row $ws4 8 7 "2024-10-26" "合肥·《四月是你的谎言》—“公生”与“薰”的钢琴小提琴唯美经典音乐集" `
    "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院" `
    "2024.10.26 19:30-10.26 21:00" 25 80 `
    "https://show.bilibili.com/platform/detail.html?id=90322" `
    "//i2.hdslb.com/bfs/openplatform/202408/BiVgXUKH1722824304648.jpeg"

Set-Row $ws4 9 8 "2024-11-09" "合肥·一生必听的钢琴曲—“从巴赫 · 莫扎特到肖邦 · 李斯特”钢琴圣手谭小棠独奏音乐会" `
    "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院" `
    "2024.11.09 19:30-11.09 21:00" 1 56 `
    "https://show.bilibili.com/platform/detail.html?id=90593" `
    "//i2.hdslb.com/bfs/openplatform/202408/SYfLxnO21723442234232.jpeg"

# Rows 10-11 no longer exist; delete the now-unused trailing rows entirely.
$ws4.Rows("10:11").Delete()
